$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.458088755607605
$ws.Range("B1").Value = 1.444346189498901
$ws.Range("C1").Value = 3.604047060012817
$ws.Range("D1").Value = 2.470268487930298
$ws.Range("E1").Value = 0.8562158942222595
